# Adjust 'Count' parameter in NBHF detections
# The header/label in column F was "Parameter 1" / data cell "nClicks" -
# both are replaced with the single, consistent label "Count" used by the
# updated source map / worksheet generator.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the column F header and its example/default value to "Count"
$ws.Range("F1").Value = "Count"
$ws.Range("F2").Value = "Count"

# The sheet used to reserve a block of empty, styled cells (G1:K1) that are
# no longer needed - clear them out so the sheet's real extent shrinks back
# down to column F.
$ws.Range("G:L").Select()
$ws.Range("G1:K1").Clear()

# Column E ("Granularity") now becomes the effective best-fit width for the
# trailing edge of the table now that the stray columns are gone.
$ws.Columns.Item(5).ColumnWidth = 9.17
